$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 2;  B = 1616.739107633333; D = 120.6574074666667 },
    @{ Row = 3;  B = 1507.6620386;      D = 111.7394087333333 },
    @{ Row = 4;  B = 1595.913612216666; D = 114.86206215 },
    @{ Row = 5;  B = 1545.856946;       D = 115.4490732666667 },
    @{ Row = 6;  B = 1605.0634982;      D = 110.6163707333333 },
    @{ Row = 7;  B = 1547.514503649999; D = 116.4257378333333 },
    @{ Row = 8;  B = 1608.401537199999; D = 118.6550773333333 },
    @{ Row = 9;  B = 1598.758053183333; D = 116.83039675 },
    @{ Row = 10; B = 1549.886868583333; D = 105.2393646833333 },
    @{ Row = 11; B = 1610.85472625;     D = 117.1247235 },
    @{ Row = 12; B = 1551.372058616666; D = 123.3677795833333 },
    @{ Row = 13; B = 1548.880400366666; D = 116.8627440333333 }
)

foreach ($entry in $values) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 4).Value = $entry.D
}
